# Applies scheduled-runner profit recalculations to the Anima Profits workbook.
# For each affected leve row, updates the price/profit columns (H-N) to the
# freshly recalculated values, including adding/removing cells where a row
# gained or lost its profit (N) figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1180.2222
$ws.Range("I58").Value = 1374.5714
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 4123.7142
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -3973.7142
$ws.Range("N58").Value = -1800
$ws.Range("H100").Value = 2234.6924
$ws.Range("J100").Value = 2327.7778
$ws.Range("L100").Value = 2327.7778
$ws.Range("N100").Value = -3409.7778
$ws.Range("H106").Value = 100006510
$ws.Range("I106").Value = 120007310
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 120007310
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -120006679
$ws.Range("N106").Value = -3762
$ws.Range("H116").Value = 8242.315000000001
$ws.Range("J116").Value = 3499.875
$ws.Range("L116").Value = 3499.875
$ws.Range("N116").Value = -10383.875
$ws.Range("H137").Value = 2463.7896
$ws.Range("I137").Value = 4375.5
$ws.Range("J137").Value = 1954
$ws.Range("K137").Value = 13126.5
$ws.Range("L137").Value = 5862
$ws.Range("M137").Value = -10576.5
$ws.Range("N137").Value = -10962
$ws.Range("H138").Value = 2179.12
$ws.Range("I138").Value = 824.25806
$ws.Range("J138").Value = 2787.8262
$ws.Range("K138").Value = 2472.77418
$ws.Range("L138").Value = 8363.4786
$ws.Range("M138").Value = 2667.22582
$ws.Range("N138").Value = -18643.4786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2940
$ws.Range("I45").Value = 2900
$ws.Range("J45").Value = 3033.3333
$ws.Range("K45").Value = 2900
$ws.Range("L45").Value = 3033.3333
$ws.Range("M45").Value = -2523
$ws.Range("N45").Value = -3787.3333
$ws.Range("H61").Value = 14496420
$ws.Range("I61").Value = 47622680
$ws.Range("J61").Value = 3681.125
$ws.Range("K61").Value = 47622680
$ws.Range("L61").Value = 3681.125
$ws.Range("M61").Value = -47622468
$ws.Range("N61").Value = -4105.125
$ws.Range("H63").Value = 38212.223
$ws.Range("I63").Value = 155350.83
$ws.Range("J63").Value = 4744.048
$ws.Range("K63").Value = 155350.83
$ws.Range("L63").Value = 4744.048
$ws.Range("M63").Value = -154664.83
$ws.Range("N63").Value = -6116.048
$ws.Range("H66").Value = 38212.223
$ws.Range("I66").Value = 155350.83
$ws.Range("J66").Value = 4744.048
$ws.Range("K66").Value = 776754.1499999999
$ws.Range("L66").Value = 23720.24
$ws.Range("M66").Value = -773322.1499999999
$ws.Range("N66").Value = -30584.24
$ws.Range("H97").Value = 1024.45
$ws.Range("I97").Value = 1535.7142
$ws.Range("J97").Value = 749.1539
$ws.Range("K97").Value = 1535.7142
$ws.Range("L97").Value = 749.1539
$ws.Range("M97").Value = -1039.7142
$ws.Range("N97").Value = -1741.1539
$ws.Range("H98").Value = 41463.75
$ws.Range("J98").Value = 41463.75
$ws.Range("L98").Value = 41463.75
$ws.Range("N98").Value = -47453.75
$ws.Range("H132").Value = 3854549.5
$ws.Range("I132").Value = 10482.2
$ws.Range("J132").Value = 7698617
$ws.Range("K132").Value = 31446.6
$ws.Range("L132").Value = 23095851
$ws.Range("M132").Value = -28916.6
$ws.Range("N132").Value = -23100911
$ws.Range("H136").Value = 14496420
$ws.Range("I136").Value = 47622680
$ws.Range("J136").Value = 3681.125
$ws.Range("K136").Value = 142868040
$ws.Range("L136").Value = 11043.375
$ws.Range("M136").Value = -142865490
$ws.Range("N136").Value = -16143.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1814.2858
$ws.Range("I99").Value = 1866.6666
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1866.6666
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -368.6666
$ws.Range("N99").Value = -4496
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H134").Value = 2903.4644
$ws.Range("I134").Value = 2535.625
$ws.Range("J134").Value = 3393.9167
$ws.Range("K134").Value = 7606.875
$ws.Range("L134").Value = 10181.7501
$ws.Range("M134").Value = -5071.875
$ws.Range("N134").Value = -15251.7501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 67500
$ws.Range("J28").Value = 67500
$ws.Range("L28").Value = 67500
$ws.Range("N28").Value = -67990
$ws.Range("H58").Value = 3354.5
$ws.Range("I58").Value = 4534.6665
$ws.Range("J58").Value = 2388.9092
$ws.Range("K58").Value = 4534.6665
$ws.Range("L58").Value = 2388.9092
$ws.Range("M58").Value = -4331.6665
$ws.Range("N58").Value = -2794.9092
$ws.Range("H122").Value = 1988.931
$ws.Range("I122").Value = 1996.7778
$ws.Range("J122").Value = 1985.4
$ws.Range("K122").Value = 5990.3334
$ws.Range("L122").Value = 5956.200000000001
$ws.Range("M122").Value = -3540.3334
$ws.Range("N122").Value = -10856.2
$ws.Range("H134").Value = 7147402
$ws.Range("I134").Value = 10875364
$ws.Range("J134").Value = 2141.8333
$ws.Range("K134").Value = 32626092
$ws.Range("L134").Value = 6425.499899999999
$ws.Range("M134").Value = -32623557
$ws.Range("N134").Value = -11495.4999
$ws.Range("H136").Value = 3354.5
$ws.Range("I136").Value = 4534.6665
$ws.Range("J136").Value = 2388.9092
$ws.Range("K136").Value = 13603.9995
$ws.Range("L136").Value = 7166.7276
$ws.Range("M136").Value = -11053.9995
$ws.Range("N136").Value = -12266.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2077.3333
$ws.Range("I132").Value = 2036.3334
$ws.Range("J132").Value = 2159.3333
$ws.Range("K132").Value = 18327.0006
$ws.Range("L132").Value = 19433.9997
$ws.Range("M132").Value = -15797.0006
$ws.Range("N132").Value = -24493.9997
$ws.Range("H136").Value = 3891.4614
$ws.Range("I136").Value = 1298.3334
$ws.Range("J136").Value = 6114.143
$ws.Range("K136").Value = 3895.0002
$ws.Range("L136").Value = 18342.429
$ws.Range("M136").Value = 1204.9998
$ws.Range("N136").Value = -28542.429
$ws.Range("H139").Value = 4006.7954
$ws.Range("I139").Value = 2104.2856
$ws.Range("J139").Value = 4894.6333
$ws.Range("K139").Value = 6312.8568
$ws.Range("L139").Value = 14683.8999
$ws.Range("M139").Value = -1172.8568
$ws.Range("N139").Value = -24963.8999
$ws.Range("H140").Value = 1749.4117
$ws.Range("I140").Value = 1316
$ws.Range("K140").Value = 3948
$ws.Range("M140").Value = 1232

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1640.6666
$ws.Range("I97").Value = 1754.6154
$ws.Range("K97").Value = 1754.6154
$ws.Range("M97").Value = -1258.6154
$ws.Range("H132").Value = 7217.304
$ws.Range("I132").Value = 13846
$ws.Range("J132").Value = 2956
$ws.Range("K132").Value = 41538
$ws.Range("L132").Value = 8868
$ws.Range("M132").Value = -39008
$ws.Range("N132").Value = -13928

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4921
$ws.Range("J100").Value = 4862
$ws.Range("L100").Value = 4862
$ws.Range("N100").Value = -5944
$ws.Range("H132").Value = 3274.12
$ws.Range("I132").Value = 2757.0625
$ws.Range("K132").Value = 8271.1875
$ws.Range("M132").Value = -5741.1875
